# copper energy redo work
# Applies the "copper energy redo work" edit:
#  1. On "pyrovshydro": move the citation text from A1 to G1 (new column),
#     re-aligning it to vertical-center (no more left-indent).
#  2. On "pyrovshydro": append forecast year rows (2022-2050) in column A.
#  3. Reposition/nudge the existing chart on "pyrovshydro" down and slightly right.
#  4. Update sheet selections (view state) to match the saved workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Sheet "pyrovshydro": move the Factbook citation from A1 -> G1
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("pyrovshydro")

$citation = $ws3.Range("A1").Value2
$ws3.Range("A1").Clear()
$ws3.Range("G1").Value = $citation
$ws3.Range("G1").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 2) Sheet "pyrovshydro": add forecast year rows 2022-2050 in column A
# ---------------------------------------------------------------------------
$startYear = 2022
$endYear = 2050
$row = 30
for ($yr = $startYear; $yr -le $endYear; $yr++) {
    $cell = $ws3.Cells.Item($row, 1)
    $cell.Value = $yr
    $cell.NumberFormat = "0"
    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 3) Sheet "pyrovshydro": nudge the existing chart's position
# ---------------------------------------------------------------------------
$co = $ws3.ChartObjects().Item(1)
$co.Top = 31.25
$co.Left = 424.265625
$co.Width = 589.625
$co.Height = 366.2499212598425

# ---------------------------------------------------------------------------
# 4) View-state bookkeeping (selections) to mirror the saved workbook
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("oregradeovertime")
$ws2.Activate()
$ws2.Range("F40").Select()

$ws3.Activate()
$ws3.Range("G1").Select()

Write-Host "copper energy redo work: edit applied"
